$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet has, for each data row (2-15), columns C..AC filled with a
# repeated "groupN" label and column AD holding a trailing count value
# (text, sometimes empty). The edit extends the "groupN" fill five more
# columns to the right (through AH) and moves the old AD value out to the
# new last column AI.

for ($row = 2; $row -le 15; $row++) {
    $styleSrc = $ws.Cells.Item($row, 29)   # column AC already carries the row's data style
    $srcCell  = $ws.Cells.Item($row, 3)    # column C holds the group label (or is blank for row 10)
    $oldLast  = $ws.Cells.Item($row, 30)   # column AD holds the old trailing value (count or blank)

    $groupVal = $srcCell.Value()
    $lastVal  = $oldLast.Value()

    # Fill the five newly-inserted columns AD..AH (30..34) with the same
    # group label as the rest of the row, matching the existing style.
    $fillRange = $ws.Range($ws.Cells.Item($row, 30), $ws.Cells.Item($row, 34))
    $fillRange.Style = $styleSrc.Style
    if ($groupVal -eq $null) {
        $fillRange.Value = ""
    } else {
        $fillRange.Value = $groupVal
    }

    # Move the old AD value into the new last column AI (35), matching style.
    # Format it as text first so a numeric-looking value (e.g. "6") is not
    # reinterpreted as a number - it was stored as text in the source column.
    $destCell = $ws.Cells.Item($row, 35)
    $destCell.Style = $styleSrc.Style
    $destCell.NumberFormat = "@"
    if ($lastVal -eq $null) {
        $destCell.Value = ""
    } else {
        $destCell.Value = $lastVal
    }
}
